$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.704868316650391
$ws.Range("B1").Value = 1.972504854202271
$ws.Range("C1").Value = 5.156136035919189
$ws.Range("D1").Value = 1.331469655036926
$ws.Range("E1").Value = 0.6571672558784485
